$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2358.674464
$ws.Range("D2").Value = 4.980593
$ws.Range("E2").Value = 0.0263

$ws.Range("B3").Value = 1475.216787
$ws.Range("D3").Value = 1.557539
$ws.Range("E3").Value = 0.212189

$ws.Range("B4").Value = 157226.247755
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = 5.715498
$ws.Range("H5").Value = -2.174512
$ws.Range("I5").Value = 13.605507
$ws.Range("J5").Value = 0.204647

$ws.Range("G6").Value = 5.363755
$ws.Range("H6").Value = -2.92343
$ws.Range("I6").Value = 13.65094
$ws.Range("J6").Value = 0.281067

$ws.Range("G7").Value = -0.351743
$ws.Range("H7").Value = -6.583153
$ws.Range("I7").Value = 5.879667
$ws.Range("J7").Value = 0.9903110000000001
